$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Sheet1: insert a new "admin/admin" credential row after the existing
# admin row, and append another "admin/admin" row at the bottom. This
# shifts the user1/user2 rows (and their hyperlinks) down by one row.
# ---------------------------------------------------------------------

# Hyperlinks don't auto-shift when rows are inserted, so drop them first
# and recreate them afterwards at the correct (shifted) addresses.
$ws1.Hyperlinks.Delete()

# Remember the plain (non-hyperlink) style used in column A so we can
# restore it after Hyperlinks.Add() re-styles the target cells.
$plainStyle = $ws1.Range("A2").Style

# Insert the new row at 3, pushing the old rows 3-5 down to 4-6.
$ws1.Rows.Item(3).Insert()
$ws1.Rows.Item(3).RowHeight = 15
$ws1.Range("A3").Value = "admin@yourstore.com"
$ws1.Range("B3").Value = "admin"

# Append a new row with the same credentials at the end (row 7).
$ws1.Range("A7").Value = "admin@yourstore.com"
$ws1.Range("A7").Style = $plainStyle
$ws1.Range("B7").Value = "admin"

# Recreate the hyperlinks on the rows that shifted down.
$ws1.Hyperlinks.Add($ws1.Range("A5"), "mailto:user1@test.com")
$ws1.Range("A5").Style = $plainStyle
$ws1.Hyperlinks.Add($ws1.Range("A6"), "mailto:user2@test.com")
$ws1.Range("A6").Style = $plainStyle

# Update the sheet's selection to match the new editing location.
$ws1.Range("A4:B7").Select()

# ---------------------------------------------------------------------
# Sheet2 (new): holds the invalid-credential rows that used to live on
# Sheet1, plus a fresh admin/admin row appended at the end.
# ---------------------------------------------------------------------

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "username"
$ws2.Range("B1").Value = "password"

$ws2.Range("A2").Value = "admin@yourstore.com"
$ws2.Range("B2").Value = "pass124"

$ws2.Range("A3").Value = "user1@test.com"
$ws2.Range("B3").Value = "admin"

$ws2.Range("A4").Value = "user2@test.com"
$ws2.Range("B4").Value = "pass126"

$ws2.Range("A5").Value = "admin@yourstore.com"
$ws2.Range("B5").Value = "admin"

$ws2.Hyperlinks.Add($ws2.Range("A3"), "mailto:user1@test.com")
$ws2.Hyperlinks.Add($ws2.Range("A4"), "mailto:user2@test.com")

# Column A on Sheet2 uses the same "credential" styling as Sheet1.
$ws2.Range("A2").Style = $plainStyle
$ws2.Range("A3").Style = $plainStyle
$ws2.Range("A4").Style = $plainStyle
$ws2.Range("A5").Style = $plainStyle

$ws2.Columns.Item(1).ColumnWidth = 19.1667

$ws2.Range("E8").Select()

# Sheet2 becomes the active/visible tab.
$ws2.Activate()
